# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Reorders / updates the worker "Periodo Mora" detail rows (16-30) on
# Hoja1 so that:
#   - the single JOSE ALBERTO PEREZ PERTUZ (doc 78031973) record now sits
#     at period 2202 (row 19) instead of row 16
#   - the MARIA CAMILA OSORIO VARGAS (doc 1052093103) records are
#     re-sequenced in ascending period order (2111..2212) and get the
#     refreshed "Valor Mora"/"Salario Basico" figures from the updated
#     database (877803 instead of 908526), except for the newest period
#     (2212) which keeps a partial "Valor Mora" of 22238.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$rows = @(
    @{ Row = 16; Tipo = "CC"; Doc = "1052093103"; Nombre = "MARIA CAMILA OSORIO VARGAS"; Periodo = "2111"; Valor = 35112;  Salario = 877803 },
    @{ Row = 17; Tipo = "CC"; Doc = "1052093103"; Nombre = "MARIA CAMILA OSORIO VARGAS"; Periodo = "2112"; Valor = 35112;  Salario = 877803 },
    @{ Row = 18; Tipo = "CC"; Doc = "1052093103"; Nombre = "MARIA CAMILA OSORIO VARGAS"; Periodo = "2201"; Valor = 35112;  Salario = 877803 },
    @{ Row = 19; Tipo = "CC"; Doc = "78031973";   Nombre = "JOSE ALBERTO PEREZ PERTUZ";  Periodo = "2202"; Valor = 40000;  Salario = 1000000 },
    @{ Row = 20; Tipo = "CC"; Doc = "1052093103"; Nombre = "MARIA CAMILA OSORIO VARGAS"; Periodo = "2202"; Valor = 35112;  Salario = 877803 },
    @{ Row = 21; Tipo = "CC"; Doc = "1052093103"; Nombre = "MARIA CAMILA OSORIO VARGAS"; Periodo = "2203"; Valor = 35112;  Salario = 877803 },
    @{ Row = 22; Tipo = "CC"; Doc = "1052093103"; Nombre = "MARIA CAMILA OSORIO VARGAS"; Periodo = "2204"; Valor = 35112;  Salario = 877803 },
    @{ Row = 23; Tipo = "CC"; Doc = "1052093103"; Nombre = "MARIA CAMILA OSORIO VARGAS"; Periodo = "2205"; Valor = 35112;  Salario = 877803 },
    @{ Row = 24; Tipo = "CC"; Doc = "1052093103"; Nombre = "MARIA CAMILA OSORIO VARGAS"; Periodo = "2206"; Valor = 35112;  Salario = 877803 },
    @{ Row = 25; Tipo = "CC"; Doc = "1052093103"; Nombre = "MARIA CAMILA OSORIO VARGAS"; Periodo = "2207"; Valor = 35112;  Salario = 877803 },
    @{ Row = 26; Tipo = "CC"; Doc = "1052093103"; Nombre = "MARIA CAMILA OSORIO VARGAS"; Periodo = "2208"; Valor = 35112;  Salario = 877803 },
    @{ Row = 27; Tipo = "CC"; Doc = "1052093103"; Nombre = "MARIA CAMILA OSORIO VARGAS"; Periodo = "2209"; Valor = 35112;  Salario = 877803 },
    @{ Row = 28; Tipo = "CC"; Doc = "1052093103"; Nombre = "MARIA CAMILA OSORIO VARGAS"; Periodo = "2210"; Valor = 35112;  Salario = 877803 },
    @{ Row = 29; Tipo = "CC"; Doc = "1052093103"; Nombre = "MARIA CAMILA OSORIO VARGAS"; Periodo = "2211"; Valor = 35112;  Salario = 877803 },
    @{ Row = 30; Tipo = "CC"; Doc = "1052093103"; Nombre = "MARIA CAMILA OSORIO VARGAS"; Periodo = "2212"; Valor = 22238;  Salario = 877803 }
)

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Cells.Item($r, 2).Value = $item.Tipo
    $ws.Cells.Item($r, 3).Value = $item.Doc
    $ws.Cells.Item($r, 4).Value = $item.Nombre
    $ws.Cells.Item($r, 5).Value = $item.Periodo
    $ws.Cells.Item($r, 6).Value = $item.Valor
    $ws.Cells.Item($r, 7).Value = $item.Salario
}
